# Fixed update to excel issue
# 1) Rename the "Requested quantity" headers on the existing sheets.
# 2) Add a new "PO Forecast" sheet (ds / PO_Forecast / yhat_lower / yhat_upper)
#    after "Monthly Trend", populated with the forecast data.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# New sheet, positioned right after "Monthly Trend"
$wsForecast = $wb.Worksheets.Add($null, $wsMonthly)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Forecast rows: ds, PO_Forecast, yhat_lower, yhat_upper
$data = @(
    @(45123.99999999999, 16, 9.12962719059146, 22.21864566924125),
    @(45130.99999999999, 16, 9.383043454587796, 22.97509222256699),
    @(45305.99999999999, 10, 3.546015680326055, 16.35350806162026),
    @(45333.99999999999, 9, 2.676465739611236, 15.68449883920265),
    @(45354.99999999999, 8, 1.867480309804237, 14.94802828828271),
    @(45361.99999999999, 8, 0.9653938201852564, 14.98410944887475),
    @(45368.99999999999, 8, 1.433291176711614, 14.86482729486391),
    @(45396.99999999999, 7, 0.7483708219831098, 13.72614931320381),
    @(45410.99999999999, 6, -0.3835598386543985, 12.59764376173176),
    @(45417.99999999999, 6, -0.691999995051094, 12.97395218126017),
    @(45424.99999999999, 6, -0.7873826441047792, 12.52665219822832),
    @(45431.99999999999, 6, -0.3916888807960443, 12.75004942992916),
    @(45438.99999999999, 5, -1.212338940578789, 12.23479297786222),
    @(45445.99999999999, 5, -1.393740032018524, 11.66376934431199),
    @(45508.99999999999, 3, -3.42563805856983, 10.49677812243933),
    @(45515.99999999999, 3, -3.430887466845317, 9.589121208867454),
    @(45536.99999999999, 2, -4.413675324378421, 8.960495599177872),
    @(45543.99999999999, 2, -4.32429116181092, 8.478464465824882),
    @(45634.99999999999, 0, -7.559937621876601, 5.096206596593663),
    @(45641.99999999999, 0, -8.457119748675607, 5.41313828204816),
    @(45648.99999999999, 0, -8.208978728027555, 4.833428558350952),
    @(45655.99999999999, 0, -8.191131866093768, 4.597097488136352),
    @(45662.99999999999, 0, -8.776894284605614, 4.30834035782599),
    @(45669.99999999999, 0, -8.97784141574739, 4.537732214343606),
    @(45676.99999999999, 0, -9.026300813279345, 4.515151305788692),
    @(45683.99999999999, 0, -9.112094714213164, 3.655002724884328),
    @(45690.99999999999, 0, -9.964426696951113, 3.85654569311634)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $wsForecast.Cells.Item($r, $c).Value = $val
        $c++
    }
    $r++
}

# Match existing formatting: bold/bordered/centered header row, date-formatted
# first column — reuse the styles already used on the other sheets instead of
# defining new ones.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A28").PasteSpecial(-4122)

$excel.CutCopyMode = $false
